{"js": "// The \"Bibliografia\" paragraph holds one long run of unbroken text.\n// Split it onto multiple lines by inserting manual line breaks\n// (Word.BreakType.line, i.e. <w:br/>) after specific anchor phrases,\n// turning:\n//   \"Bibliografia b\u00e1sica:PRESS...656p.REED...508p.Bibliografia\n//    complementar:TEIXEIRA...623p.\"\n// into 5 lines (with a blank line between the two bibliography\n// sections), all still inside the same run.\n\n// Anchors are searched for, in document order, and a line break is\n// inserted right after each one. \"508p.\" (the one that precedes\n// \"Bibliografia complementar:\") gets two consecutive breaks to create\n// the blank line seen in the target markup.\nconst anchors = [\n  { text: \"b\u00e1sica:\", breaks: 1 },\n  { text: \"656p.\", breaks: 1 },\n  { text: \"508p.\", breaks: 2 },\n  { text: \"complementar:\", breaks: 1 },\n];\n\nfor (const anchor of anchors) {\n  const results = context.document.body.search(anchor.text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  const found = results.items[0];\n  // Office.js represents a manual line break (<w:br/>) as the\n  // vertical-tab control character (U+000B) inside run text.\n  found.insertText(\"\\u000b\".repeat(anchor.breaks), Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The bibliography paragraph currently holds one long run of text with no\n# line breaks. Split it into several lines (using manual line breaks, i.e.\n# <w:br/>) while keeping everything inside the same run, by doing a series\n# of targeted Find/Replace operations using the \"^l\" (manual line break)\n# replacement code.\n$replacements = @(\n    @(\"b\u00e1sica:PRESS\", \"b\u00e1sica:^lPRESS\"),\n    @(\"656p.REED\", \"656p.^lREED\"),\n    @(\"508p.Bibliografia complementar:\", \"508p.^l^lBibliografia complementar:\"),\n    @(\"complementar:TEIXEIRA\", \"complementar:^lTEIXEIRA\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $null = $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
